$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Updated cryptos list" GitHub Actions refresh: new price/volume figures
# for every coin row, plus three pairs of adjacent rows (14/15, 42/43,
# 48/49) that swapped rank order as their prices moved.

# Column D ("Price") holds scraped text, not numbers -- some values look
# like plain decimals ("1.000", "0.4980") that Excel would otherwise
# auto-convert and trim (1 instead of 1.000). Force Text format on every
# Price cell we are about to write so the literal digits are preserved,
# matching how the rest of the column is already stored.
$forceTextCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '26.895.48'
$ws.Range("E2").Value = '  -1.74%  '

# Row 3
$ws.Range("D3").Value = '1.808.57'
$ws.Range("E3").Value = '  -0.98%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '310.39'
$ws.Range("E5").Value = '  -1.04%  '

# Row 6
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("D7").Value = '0.4617'
$ws.Range("E7").Value = '  +3.44%  '

# Row 8
$ws.Range("D8").Value = '0.3707'
$ws.Range("E8").Value = '  -1.98%  '

# Row 9
$ws.Range("D9").Value = '0.07387'
$ws.Range("E9").Value = '  -0.31%  '

# Row 10
$ws.Range("D10").Value = '0.8743'
$ws.Range("E10").Value = '  -0.69%  '

# Row 11
$ws.Range("E11").Value = '  -1.99%  '

# Row 12
$ws.Range("D12").Value = '1.812.14'
$ws.Range("E12").Value = '  -0.85%  '

# Row 13
$ws.Range("D13").Value = '5.359'
$ws.Range("E13").Value = '  -1.44%  '

# Row 14
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '92.09'
$ws.Range("E14").Value = '  -0.66%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '6.519'
$ws.Range("E15").Value = '  -3.00%  '

# Row 16
$ws.Range("D16").Value = '0.07043'
$ws.Range("E16").Value = '  -0.23%  '

# Row 17
$ws.Range("E17").Value = '  +0.06%  '

# Row 18
$ws.Range("D18").Value = '0.000008693'
$ws.Range("E18").Value = '  -1.34%  '

# Row 19
$ws.Range("E19").Value = '  +0.10%  '

# Row 20
$ws.Range("D20").Value = '14.75'
$ws.Range("E20").Value = '  -2.14%  '

# Row 21
$ws.Range("D21").Value = '26.890.45'
$ws.Range("E21").Value = '  -1.80%  '

# Row 22
$ws.Range("D22").Value = '5.326'
$ws.Range("E22").Value = '  -0.78%  '

# Row 23
$ws.Range("D23").Value = '10.65'
$ws.Range("E23").Value = '  -2.84%  '

# Row 24
$ws.Range("D24").Value = '2.021.10'
$ws.Range("E24").Value = '  -1.44%  '

# Row 25
$ws.Range("D25").Value = '1.893'
$ws.Range("E25").Value = '  -3.47%  '

# Row 26
$ws.Range("D26").Value = '151.33'
$ws.Range("E26").Value = '  +0.26%  '

# Row 27
$ws.Range("D27").Value = '18.35'
$ws.Range("E27").Value = '  -1.40%  '

# Row 28
$ws.Range("D28").Value = '2.148'
$ws.Range("E28").Value = '  -6.03%  '

# Row 29
$ws.Range("D29").Value = '5.313'
$ws.Range("E29").Value = '  -0.92%  '

# Row 30
$ws.Range("E30").Value = '  -1.04%  '

# Row 31
$ws.Range("D31").Value = '0.08901'
$ws.Range("E31").Value = '  -0.06%  '

# Row 32
$ws.Range("D32").Value = '0.7531'
$ws.Range("E32").Value = '  -5.10%  '

# Row 33
$ws.Range("D33").Value = '1.159'
$ws.Range("E33").Value = '  -3.19%  '

# Row 34
$ws.Range("D34").Value = '4.445'
$ws.Range("E34").Value = '  -2.88%  '

# Row 35
$ws.Range("D35").Value = '2.910'
$ws.Range("E35").Value = '  -0.67%  '

# Row 36
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  +0.04%  '

# Row 37
$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  -0.51%  '

# Row 38
$ws.Range("D38").Value = '0.01970'
$ws.Range("E38").Value = '  -0.82%  '

# Row 39
$ws.Range("D39").Value = '0.05245'
$ws.Range("E39").Value = '  -0.60%  '

# Row 40
$ws.Range("D40").Value = '2.423'
$ws.Range("E40").Value = '  +2.17%  '

# Row 41
$ws.Range("D41").Value = '2.928'
$ws.Range("E41").Value = '  +1.93%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '7.226'
$ws.Range("E42").Value = '  -1.33%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.5325'
$ws.Range("E43").Value = '  +0.18%  '

# Row 44
$ws.Range("D44").Value = '0.1662'
$ws.Range("E44").Value = '  -2.24%  '

# Row 45
$ws.Range("D45").Value = '8.521'
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("D46").Value = '0.4980'
$ws.Range("E46").Value = '  -1.54%  '

# Row 47
$ws.Range("D47").Value = '10.32'
$ws.Range("E47").Value = '  -2.95%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.672'
$ws.Range("E48").Value = '  -1.02%  '

# Row 49
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  +0.06%  '

# Row 50
$ws.Range("D50").Value = '103.72'
$ws.Range("E50").Value = '  -1.83%  '

# Row 51
$ws.Range("D51").Value = '0.06293'
$ws.Range("E51").Value = '  -1.49%  '
